$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.42"
$ws.Range("E2").Value = "'2.32%"
$ws.Range("D3").Value = "'35.20"
$ws.Range("E3").Value = "'13.32%"
$ws.Range("D4").Value = "'5.146"
$ws.Range("E4").Value = "'4.87%"
$ws.Range("D5").Value = "'0.07771"
$ws.Range("E5").Value = "'4.44%"
$ws.Range("E6").Value = "'7.25%"
$ws.Range("D7").Value = "'8.039"
$ws.Range("E7").Value = "'3.66%"
$ws.Range("D8").Value = "'3.944"
$ws.Range("E8").Value = "'5.22%"
$ws.Range("D9").Value = "'0.9323"
$ws.Range("E9").Value = "'1.90%"
$ws.Range("D10").Value = "'0.09926"
$ws.Range("E10").Value = "'11.60%"
$ws.Range("D11").Value = "'0.1795"
$ws.Range("E11").Value = "'4.89%"
$ws.Range("D12").Value = "'0.08605"
$ws.Range("E12").Value = "'3.13%"
$ws.Range("D13").Value = "'0.03320"
$ws.Range("E13").Value = "'5.11%"
$ws.Range("D14").Value = "'0.09917"
$ws.Range("E14").Value = "'-1.72%"
$ws.Range("D15").Value = "'0.001499"
$ws.Range("E15").Value = "'-1.59%"
$ws.Range("D16").Value = "'0.005855"
$ws.Range("E16").Value = "'1.43%"
$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'-1.41%"
$ws.Range("D19").Value = "'0.3367"
$ws.Range("E19").Value = "'1.13%"
$ws.Range("D21").Value = "'4.320"
$ws.Range("E21").Value = "'8.79%"
$ws.Range("D22").Value = "'0.2303"
$ws.Range("E22").Value = "'9.54%"
$ws.Range("D23").Value = "'0.04560"
$ws.Range("E23").Value = "'0.10%"
$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'0.13%"
$ws.Range("E25").Value = "'-5.35%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'0.09%"
$ws.Range("E27").Value = "'-0.16%"
$ws.Range("D39").Value = "'0.01779"
$ws.Range("E39").Value = "'10.06%"
$ws.Range("D40").Value = "'0.04795"
$ws.Range("E40").Value = "'6.79%"
$ws.Range("D41").Value = "'0.007792"
$ws.Range("E41").Value = "'6.89%"
$ws.Range("E42").Value = "'6.18%"
$ws.Range("D43").Value = "'0.006814"
$ws.Range("E43").Value = "'-24.73%"
$ws.Range("D44").Value = "'0.002081"
$ws.Range("E44").Value = "'5.58%"
$ws.Range("D45").Value = "'0.009449"
$ws.Range("E45").Value = "'3.34%"
$ws.Range("D46").Value = "'0.00006114"
$ws.Range("E46").Value = "'-0.15%"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'2.962"
$ws.Range("E48").Value = "'32.50%"
$ws.Range("D49").Value = "'0.002002"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.06%"
